# Summer 23 Week 11 update: append 12 new matchup rows (1059-1070) to sheet "Nine"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newData = @(
    @(5, 8, 6, 12),
    @(3, 18, 4, 2),
    @(5, 13, 6, 7),
    @(8, 16, 7, 4),
    @(4, 15, 3, 5),
    @(4, 7, 3, 13),
    @(3, 16, 4, 4),
    @(5, 4, 3, 16),
    @(5, 4, 9, 16),
    @(4, 3, 3, 17),
    @(3, 13, 4, 7),
    @(4, 18, 3, 2)
)

$startRow = 1059
for ($i = 0; $i -lt $newData.Count; $i++) {
    $r = $startRow + $i
    $row = $newData[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Update the window/view to match scrolled-down state seen after adding rows
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 1055
$win.ScrollColumn = 1
$ws.Range("A1071").Select()
